$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "F2", "F2rl1", "ECs",  3, 1, 0.6740386666666667, 2.022116,           0.3058319194585966, 0.3058319194585966, 1, 0.3333333333333333, 0.1341563333333333, 0.402469,            0.02276964089369097, 0.02276964089369097, 0.0904265560448889,  0.813839004404,       0.006963682979900465, 0.006963682979900464),
    @("ECs",  "F2", "F2rl1", "FAPs", 3, 1, 0.6740386666666667, 2.022116,           0.3058319194585966, 0.3058319194585966, 3, 1,                  5.757738666666666,  17.273216,           0.9772303591063091,  0.9772303591063091,  3.880938493895111,   34.928446445056,      0.2988682364786961,   0.2988682364786961),
    @("FAPs", "F2", "F2rl1", "ECs",  3, 1, 0.9731926666666667, 2.919578,           0.4415672215387696, 0.4415672215387696, 1, 0.3333333333333333, 0.1341563333333333, 0.402469,            0.02276964089369097, 0.02276964089369097, 0.1305599597868889,  1.175039638082,       0.01005432706486267,  0.01005432706486267),
    @("FAPs", "F2", "F2rl1", "FAPs", 3, 1, 0.9731926666666667, 2.919578,           0.4415672215387696, 0.4415672215387696, 3, 1,                  5.757738666666666,  17.273216,           0.9772303591063091,  0.9772303591063091,  5.603389046983111,   50.430501422848,      0.4315128944739069,   0.4315128944739069),
    @("M2",   "F2", "F2rl1", "ECs",  1, 0.3333333333333333, 0.048595,            0.145785,             0.02204903496054208, 0.02204903496054208, 1, 0.3333333333333333, 0.1341563333333333, 0.402469,            0.02276964089369097, 0.02276964089369097, 0.006519327018333334, 0.058673943165,     0.0005020486081039809, 0.0005020486081039807),
    @("M2",   "F2", "F2rl1", "FAPs", 1, 0.3333333333333333, 0.048595,            0.145785,             0.02204903496054208, 0.02204903496054208, 3, 1,                  5.757738666666666,  17.273216,           0.9772303591063091,  0.9772303591063091,  0.2797973105066666,  2.51817579456,       0.0215469863524381,   0.0215469863524381),
    @("sCs",  "F2", "F2rl1", "ECs",  3, 1, 0.508125,           1.524375,           0.2305518240420917, 0.2305518240420917, 1, 0.3333333333333333, 0.1341563333333333, 0.402469,            0.02276964089369097, 0.02276964089369097, 0.06816818687500001, 0.613513681875,      0.005249582240823857, 0.005249582240823856),
    @("sCs",  "F2", "F2rl1", "FAPs", 3, 1, 0.508125,           1.524375,           0.2305518240420917, 0.2305518240420917, 3, 1,                  5.757738666666666,  17.273216,           0.9772303591063091,  0.9772303591063091,  2.92565096,          26.33085864,          0.2253022418012678,   0.2253022418012678)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}
